# Review_144 -> Review_143 edit: new title/links, new Hebrew review text,
# an inserted blank "Normal" paragraph before the body text, and the
# trailing paragraph restyled from "Normal" to "Heading 2".

$d = $word.ActiveDocument

# --- 1. Heading paragraph: new title + link, drop the stray tab before the URL ---
$pHeading = $d.Paragraphs(1)
$pHeading.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Review 143: [Short] Explaining grokking through circuit efficiency,  11.09.2023</w:t><w:br/><w:t>https://arxiv.org/abs/2309.02390</w:t></w:r></w:p>') | Out-Null

# --- 2. Bold "Paper:" paragraph: new arXiv link ---
$pPaper = $d.Paragraphs(2)
$pPaper.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:b/></w:rPr><w:t>Paper: https://arxiv.org/abs/2309.02390v1</w:t></w:r></w:p>') | Out-Null

# --- 3. Insert a new blank "Normal" paragraph right before the review body ---
$pOldBody = $d.Paragraphs(4)
$pOldBody.Range.InsertParagraphBefore() | Out-Null

# --- 4. Replace the review body paragraph (now shifted down by one) ---
$pBody = $d.Paragraphs(5)
$pBody.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:t>יש תופעה מעניינת, הנקראת גרוקינג (grokking) המתרחשת (לפעמים) במהלך אימון של רשתות נוירונים. למעשה גרוקינג מחלק אימון של רשת נוירונים לשני שלבים עיקריים: הראשון הוא שלב השינון (memorization) כאשר הרשת משננת את הדוגמאות(overfit) ושלב ההכללה (למידה אמיתית).</w:t><w:br/><w:t xml:space="preserve">היום ב-#shorthebrewpapereviews אנחנו סוקרים מאמר המנסה להסביר למה התופעה הזו מתרחשת. למה בכלל רשת נוירונים ״משננת״ את הדוגמאות בתחילת האימון מרגע מסוים עוברת למשטר ההכללה, כלומר למידה אמיתית. </w:t><w:br/><w:br/><w:t xml:space="preserve">קודם כל המחברים שמו לב שבסוף משטר השינון הדיוק של המודל הוא כמעט מושלם  והלוס מאד נמוך ולמרות זאת מרגע מסוים הרשת מתחילה ללמוד להכליל משום מה למרות זאת. המחברים טוענים שהסיבה לכך היא הרשת מצליחה להוריד את הלוס עוד יותר (במצב של דיוק כמעט מושלם על הטריין סט) על ידי כך שהיא מתחילה להוציא חיזויים בטוחים יותר (עבור הקטגוריות הנכונות). </w:t><w:br/><w:br/><w:t xml:space="preserve">אני חושב שראיתי מאמרים(לדעתי בנושא double descent) שטוענים שהירידה בלוס כזו במצב של ערכי לוס מאוד נמוכים הופכת את הפונקציה שהרשת ממדלת ליותר פשוטה (מאפס הרבה פרמטרים ומגדיל את האחרים). הם גם מתארים תופעה מעניינת נוספת: כאשר במשטר השינון מקפיאים את לוג''יטים (שמהם מחשבים את ההסתברות הקטגוריות) הנורמה של פרמטרי הרשת עולה עם הגדלת הדאטהסט. </w:t><w:br/><w:br/><w:t>לעומת זאת במשטר ההכללה נורמה של פרמטרי הרשת לא תלויה בגודל הדאטהסט. בנוסף הם שמו לב שלפעמים ניתן לצפות תופעה שהם קראו לה semi-grokking כאשר המעבר למצב שינון מתרחש אבל הרשת מצליחה להגיע למצב ההכללה אבל לא מצליחה להגיע להכללה טובה ונתקעת איפשהו באמצע (semi-overfit). בקיצור מאמר מאוד מעניין – מומלץ בחום!</w:t></w:r></w:p>') | Out-Null

# --- 5. Restyle the trailing empty paragraph: Normal -> Heading 2 ---
$pTrailing = $d.Paragraphs(6)
$pTrailing.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r/></w:p>') | Out-Null

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
